# Updated cryptos list - apply new Price/Volume(1h) text values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.010.53"
$ws.Range("E2").Value = "  -3.74%  "
$ws.Range("D3").Value = "3.160.94"
$ws.Range("E3").Value = "  -8.89%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "560.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.37"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.49%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  -3.83%  "
$ws.Range("D9").Value = "3.160.45"
$ws.Range("E9").Value = "  -8.89%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.123"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.69%  "
$ws.Range("E11").Value = "  -5.59%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.393"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.85%  "
$ws.Range("D13").Value = "3.709.40"
$ws.Range("E13").Value = "  -8.85%  "
$ws.Range("E14").Value = "  +1.31%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.11"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -10.63%  "
$ws.Range("D16").Value = "64.014.82"
$ws.Range("E16").Value = "  -3.58%  "
$ws.Range("E17").Value = "  -6.08%  "
$ws.Range("D18").Value = "3.163.70"
$ws.Range("E18").Value = "  -8.87%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.70"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.88%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.89"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.96%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "350.55"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.13"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -7.20%  "
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.42"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.81%  "
$ws.Range("E25").Value = "  -6.40%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.500"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.96%  "
$ws.Range("E27").Value = "  -5.30%  "
$ws.Range("E28").Value = "  -1.08%  "
$ws.Range("E29").Value = "  +0.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.47"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -8.58%  "
$ws.Range("E32").Value = "  -5.88%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.86"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -7.91%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.60"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.65%  "
$ws.Range("E35").Value = "  -6.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.43"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -8.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "153.06"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.813"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -8.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "25.67"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -9.43%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.68"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.20%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.48"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.92%  "
$ws.Range("D42").Value = "2.591.67"
$ws.Range("E42").Value = "  -7.15%  "
$ws.Range("E43").Value = "  -7.99%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "39.32"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.91%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.94"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -8.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0647"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.96%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "23.53"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "317.35"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -7.58%  "
$ws.Range("E49").Value = "  -8.53%  "
$ws.Range("E50").Value = "  -4.28%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.999"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.06%  "
